{"js": "// Replace the two-digit multiplication problems/answers in the table\n// with the new generated values. All 25 original values are unique in the\n// document, so a direct search+replace per pair is unambiguous.\nconst replacements = [\n  [\"38\u00d718=684\", \"42\u00d799=4158\"],\n  [\"71\u00d715=1065\", \"36\u00d756=2016\"],\n  [\"24\u00d768=1632\", \"97\u00d716=1552\"],\n  [\"98\u00d765=6370\", \"16\u00d793=1488\"],\n  [\"22\u00d764=1408\", \"89\u00d797=8633\"],\n  [\"50\u00d782=4100\", \"92\u00d711=1012\"],\n  [\"17\u00d735=595\", \"41\u00d712=492\"],\n  [\"71\u00d794=6674\", \"38\u00d773=2774\"],\n  [\"24\u00d795=2280\", \"38\u00d741=1558\"],\n  [\"53\u00d722=1166\", \"28\u00d765=1820\"],\n  [\"21\u00d743=903\", \"16\u00d772=1152\"],\n  [\"77\u00d791=7007\", \"12\u00d749=588\"],\n  [\"81\u00d725=2025\", \"93\u00d729=2697\"],\n  [\"35\u00d756=1960\", \"81\u00d795=7695\"],\n  [\"57\u00d712=684\", \"13\u00d757=741\"],\n  [\"11\u00d784=924\", \"33\u00d784=2772\"],\n  [\"30\u00d778=2340\", \"93\u00d727=2511\"],\n  [\"53\u00d787=4611\", \"56\u00d767=3752\"],\n  [\"79\u00d797=7663\", \"70\u00d724=1680\"],\n  [\"58\u00d787=5046\", \"92\u00d750=4600\"],\n  [\"32\u00d729=928\", \"96\u00d784=8064\"],\n  [\"48\u00d783=3984\", \"50\u00d738=1900\"],\n  [\"82\u00d756=4592\", \"60\u00d771=4260\"],\n  [\"73\u00d778=5694\", \"43\u00d739=1677\"],\n  [\"58\u00d791=5278\", \"18\u00d716=288\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication problems/answers in the table\n# with the new generated values. All 25 original values are unique in the\n# document, so a direct Find/Replace per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n    @(\"38\u00d718=684\", \"42\u00d799=4158\"),\n    @(\"71\u00d715=1065\", \"36\u00d756=2016\"),\n    @(\"24\u00d768=1632\", \"97\u00d716=1552\"),\n    @(\"98\u00d765=6370\", \"16\u00d793=1488\"),\n    @(\"22\u00d764=1408\", \"89\u00d797=8633\"),\n    @(\"50\u00d782=4100\", \"92\u00d711=1012\"),\n    @(\"17\u00d735=595\", \"41\u00d712=492\"),\n    @(\"71\u00d794=6674\", \"38\u00d773=2774\"),\n    @(\"24\u00d795=2280\", \"38\u00d741=1558\"),\n    @(\"53\u00d722=1166\", \"28\u00d765=1820\"),\n    @(\"21\u00d743=903\", \"16\u00d772=1152\"),\n    @(\"77\u00d791=7007\", \"12\u00d749=588\"),\n    @(\"81\u00d725=2025\", \"93\u00d729=2697\"),\n    @(\"35\u00d756=1960\", \"81\u00d795=7695\"),\n    @(\"57\u00d712=684\", \"13\u00d757=741\"),\n    @(\"11\u00d784=924\", \"33\u00d784=2772\"),\n    @(\"30\u00d778=2340\", \"93\u00d727=2511\"),\n    @(\"53\u00d787=4611\", \"56\u00d767=3752\"),\n    @(\"79\u00d797=7663\", \"70\u00d724=1680\"),\n    @(\"58\u00d787=5046\", \"92\u00d750=4600\"),\n    @(\"32\u00d729=928\", \"96\u00d784=8064\"),\n    @(\"48\u00d783=3984\", \"50\u00d738=1900\"),\n    @(\"82\u00d756=4592\", \"60\u00d771=4260\"),\n    @(\"73\u00d778=5694\", \"43\u00d739=1677\"),\n    @(\"58\u00d791=5278\", \"18\u00d716=288\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n}\n"}
